$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Append "wifi_Mode(ON);" (replacing the trailing newline) to the script
# steps for the VT200-0406 (G3) and VT200-0407 (G4) test rows.
$ws.Range("G3").Value = "wait(5);`nvalidate1;`nlink_Click(network_test_link);`nwait(5);`nvalidate2;`nSelectTestToRun(VT200_0406_string);`nClickRunTest(runtest_top_xpath);`nwait(5);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nvalidate4;`npress_Key(Home);`nvalidate5;`nwifi_Mode(OFF);`nlaunch_App_Device(com.rhomobile.compliancejs);`nwait(5);`nSelectTestToRun(VT200_0406_string);`nClickRunTest(runtest_top_xpath);`nwait(5);`nvalidate6;`nClickRunTest(runtest_bottom_xpath);`nvalidate7;`nwifi_Mode(ON);"

$ws.Range("G4").Value = "wait(5);`nvalidate1;`nlink_Click(network_test_link);`nwait(5);`nvalidate2;`nSelectTestToRun(VT200_0407_string);`nClickRunTest(runtest_top_xpath);`nwait(5);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nvalidate4;`npress_Key(Home);`nvalidate5;`nwifi_Mode(OFF);`nlaunch_App_Device(com.rhomobile.compliancejs);`nwait(5);`nSelectTestToRun(VT200_0407_string);`nClickRunTest(runtest_top_xpath);`nwait(5);`nvalidate6;`nClickRunTest(runtest_bottom_xpath);`nvalidate7;`nwifi_Mode(ON);"

# Move the active selection to G4, matching the saved view state.
$ws.Range("G4").Select()
